$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = '23.733.22'
$ws.Cells.Item(2, 5).Value = '  +1.09%  '

$ws.Cells.Item(3, 4).Value = '1.658.75'
$ws.Cells.Item(3, 5).Value = '  +1.09%  '

Set-CellText $ws.Cells.Item(4, 4) '1.000'
$ws.Cells.Item(4, 5).Value = '  -0.02%  '

Set-CellText $ws.Cells.Item(5, 4) '0.9998'
$ws.Cells.Item(5, 5).Value = '  +0.09%  '

Set-CellText $ws.Cells.Item(6, 4) '303.65'
$ws.Cells.Item(6, 5).Value = '  -0.14%  '

Set-CellText $ws.Cells.Item(7, 4) '0.3812'
$ws.Cells.Item(7, 5).Value = '  +0.55%  '

$ws.Cells.Item(8, 5).Value = '  -0.35%  '

$ws.Cells.Item(9, 5).Value = '  -0.87%  '

Set-CellText $ws.Cells.Item(10, 4) '0.08214'
$ws.Cells.Item(10, 5).Value = '  +0.15%  '

Set-CellText $ws.Cells.Item(11, 4) '1.242'
$ws.Cells.Item(11, 5).Value = '  +0.50%  '

$ws.Cells.Item(12, 5).Value = '  -0.02%  '

Set-CellText $ws.Cells.Item(13, 4) '22.69'
$ws.Cells.Item(13, 5).Value = '  +0.56%  '

Set-CellText $ws.Cells.Item(14, 4) '6.529'
$ws.Cells.Item(14, 5).Value = '  +0.89%  '

Set-CellText $ws.Cells.Item(15, 4) '7.432'
$ws.Cells.Item(15, 5).Value = '  +0.43%  '

Set-CellText $ws.Cells.Item(16, 4) '0.00001234'
$ws.Cells.Item(16, 5).Value = '  -0.70%  '

$ws.Cells.Item(17, 4).Value = '1.648.10'
$ws.Cells.Item(17, 5).Value = '  +1.01%  '

$ws.Cells.Item(18, 5).Value = '  +2.28%  '

Set-CellText $ws.Cells.Item(19, 4) '0.07002'
$ws.Cells.Item(19, 5).Value = '  +0.96%  '

Set-CellText $ws.Cells.Item(20, 4) '6.843'
$ws.Cells.Item(20, 5).Value = '  +3.85%  '

Set-CellText $ws.Cells.Item(21, 4) '17.70'
$ws.Cells.Item(21, 5).Value = '  +0.92%  '

$ws.Cells.Item(22, 5).Value = '  +0.12%  '

$ws.Cells.Item(23, 5).Value = '  +2.30%  '

$ws.Cells.Item(24, 4).Value = '23.741.23'
$ws.Cells.Item(24, 5).Value = '  +1.12%  '

Set-CellText $ws.Cells.Item(25, 4) '2.516'
$ws.Cells.Item(25, 5).Value = '  +0.32%  '

Set-CellText $ws.Cells.Item(26, 4) '3.053'
$ws.Cells.Item(26, 5).Value = '  -0.31%  '

Set-CellText $ws.Cells.Item(27, 4) '21.31'
$ws.Cells.Item(27, 5).Value = '  +0.62%  '

Set-CellText $ws.Cells.Item(28, 4) '151.95'
$ws.Cells.Item(28, 5).Value = '  +0.36%  '

Set-CellText $ws.Cells.Item(29, 4) '5.207'
$ws.Cells.Item(29, 5).Value = '  -1.23%  '

Set-CellText $ws.Cells.Item(30, 4) '134.63'
$ws.Cells.Item(30, 5).Value = '  +0.72%  '

$ws.Cells.Item(31, 4).Value = '1.838.22'
$ws.Cells.Item(31, 5).Value = '  +1.11%  '

Set-CellText $ws.Cells.Item(32, 4) '6.982'
$ws.Cells.Item(32, 5).Value = '  +4.76%  '

Set-CellText $ws.Cells.Item(33, 4) '2.189'
$ws.Cells.Item(33, 5).Value = '  +0.18%  '

Set-CellText $ws.Cells.Item(36, 4) '0.02822'
$ws.Cells.Item(36, 5).Value = '  +2.00%  '

Set-CellText $ws.Cells.Item(37, 4) '0.2528'
$ws.Cells.Item(37, 5).Value = '  +1.20%  '

Set-CellText $ws.Cells.Item(38, 4) '6.127'
$ws.Cells.Item(38, 5).Value = '  +1.63%  '

$ws.Cells.Item(39, 5).Value = '  +0.12%  '

Set-CellText $ws.Cells.Item(40, 4) '0.07074'
$ws.Cells.Item(40, 5).Value = '  -0.84%  '

Set-CellText $ws.Cells.Item(41, 4) '13.10'
$ws.Cells.Item(41, 5).Value = '  +7.67%  '

Set-CellText $ws.Cells.Item(42, 4) '0.7043'
$ws.Cells.Item(42, 5).Value = '  -0.37%  '

Set-CellText $ws.Cells.Item(43, 4) '1.338'
$ws.Cells.Item(43, 5).Value = '  -0.17%  '

Set-CellText $ws.Cells.Item(44, 4) '16.02'
$ws.Cells.Item(44, 5).Value = '  +1.22%  '

Set-CellText $ws.Cells.Item(45, 4) '0.6548'
$ws.Cells.Item(45, 5).Value = '  -0.19%  '

Set-CellText $ws.Cells.Item(46, 4) '2.329'
$ws.Cells.Item(46, 5).Value = '  +1.75%  '

Set-CellText $ws.Cells.Item(47, 4) '0.9991'
$ws.Cells.Item(47, 5).Value = '  +0.00%  '

Set-CellText $ws.Cells.Item(48, 4) '3.976'
$ws.Cells.Item(48, 5).Value = '  +0.25%  '

Set-CellText $ws.Cells.Item(49, 4) '0.07953'
$ws.Cells.Item(49, 5).Value = '  -0.47%  '

Set-CellText $ws.Cells.Item(50, 4) '128.28'
$ws.Cells.Item(50, 5).Value = '  +0.45%  '

Set-CellText $ws.Cells.Item(51, 4) '1.186'
$ws.Cells.Item(51, 5).Value = '  -0.85%  '

# Row 34: now ImmutableX
$ws.Cells.Item(34, 2).Value = 'ImmutableX'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws.Cells.Item(34, 4) '1.064'
$ws.Cells.Item(34, 5).Value = '  -0.03%  '

# Row 35: now FraxShare
$ws.Cells.Item(35, 2).Value = 'FraxShare'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws.Cells.Item(35, 4) '11.96'
$ws.Cells.Item(35, 5).Value = '  +4.77%  '
